$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H11").Value = 46.133335
$ws.Range("I11").Value = 46.133335
$ws.Range("K11").Value = 46.133335
$ws.Range("M11").Value = 93.866665
$ws.Range("H12").Value = 396.6
$ws.Range("I12").Value = 240.5
$ws.Range("K12").Value = 240.5
$ws.Range("M12").Value = -70.5
$ws.Range("H19").Value = 422.22223
$ws.Range("I19").Value = 250
$ws.Range("J19").Value = 471.42856
$ws.Range("K19").Value = 250
$ws.Range("L19").Value = 471.42856
$ws.Range("M19").Value = -75
$ws.Range("N19").Value = -821.4285600000001
$ws.Range("H69").Value = 7507.25
$ws.Range("I69").Value = 5999.5
$ws.Range("K69").Value = 17998.5
$ws.Range("M69").Value = -17124.5
$ws.Range("H72").Value = 7507.25
$ws.Range("I72").Value = 5999.5
$ws.Range("K72").Value = 53995.5
$ws.Range("M72").Value = -49627.5
$ws.Range("H88").Value = 2017.625
$ws.Range("I88").Value = 3219.8
$ws.Range("J88").Value = 1471.1818
$ws.Range("K88").Value = 3219.8
$ws.Range("L88").Value = 1471.1818
$ws.Range("M88").Value = -2813.8
$ws.Range("N88").Value = -2283.1818
$ws.Range("H91").Value = 2017.625
$ws.Range("I91").Value = 3219.8
$ws.Range("J91").Value = 1471.1818
$ws.Range("K91").Value = 3219.8
$ws.Range("L91").Value = 1471.1818
$ws.Range("M91").Value = -1815.8
$ws.Range("N91").Value = -4279.1818
$ws.Range("H137").Value = 1414.1428
$ws.Range("J137").Value = 1339.6
$ws.Range("L137").Value = 4018.8
$ws.Range("N137").Value = -9118.799999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 1910.6666
$ws.Range("I61").Value = 1910.6666
$ws.Range("K61").Value = 1910.6666
$ws.Range("M61").Value = -1698.6666
$ws.Range("H102").Value = 2636.3333
$ws.Range("I102").Value = 2636.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2636.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1014.3333
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 4887
$ws.Range("I122").Value = 4887
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14661
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12211
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 1910.6666
$ws.Range("I136").Value = 1910.6666
$ws.Range("K136").Value = 5731.9998
$ws.Range("M136").Value = -3181.9998

# --- Sheet BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Range("H20").Value = 1486.9286
$ws.Range("I20").Value = 1414.2222
$ws.Range("J20").Value = 1617.8
$ws.Range("K20").Value = 1414.2222
$ws.Range("L20").Value = 1617.8
$ws.Range("M20").Value = -1167.2222
$ws.Range("N20").Value = -2111.8
$ws.Range("H86").Value = 30991.4
$ws.Range("J86").Value = 50669
$ws.Range("L86").Value = 50669
$ws.Range("N86").Value = -52915
$ws.Range("H89").Value = 30991.4
$ws.Range("J89").Value = 50669
$ws.Range("L89").Value = 253345
$ws.Range("N89").Value = -264577
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H117").Value = 99900
$ws.Range("J117").Value = 99900
$ws.Range("L117").Value = 99900
$ws.Range("N117").Value = -109078

# --- Sheet CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 8
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 166
$ws.Range("N17").ClearContents()
$ws.Range("H31").Value = 1945.0834
$ws.Range("I31").Value = 1231.625
$ws.Range("K31").Value = 1231.625
$ws.Range("M31").Value = -936.625
$ws.Range("H34").Value = 1945.0834
$ws.Range("I34").Value = 1231.625
$ws.Range("K34").Value = 1231.625
$ws.Range("M34").Value = -1029.625
$ws.Range("H41").Value = 9514.75
$ws.Range("I41").Value = 6059
$ws.Range("J41").Value = 10666.667
$ws.Range("K41").Value = 6059
$ws.Range("L41").Value = 10666.667
$ws.Range("M41").Value = -5631
$ws.Range("N41").Value = -11522.667
$ws.Range("H50").Value = 23714.285
$ws.Range("I50").Value = 24000
$ws.Range("J50").Value = 23500
$ws.Range("K50").Value = 24000
$ws.Range("L50").Value = 23500
$ws.Range("M50").Value = -23375
$ws.Range("N50").Value = -24750
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 11337.5
$ws.Range("I60").Value = 6718.1816
$ws.Range("J60").Value = 21500
$ws.Range("K60").Value = 6718.1816
$ws.Range("L60").Value = 21500
$ws.Range("M60").Value = -6207.1816
$ws.Range("N60").Value = -22522
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 6496
$ws.Range("I132").Value = 6496
$ws.Range("K132").Value = 19488
$ws.Range("M132").Value = -16958

# --- Sheet CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H14").Value = 507.42856
$ws.Range("I14").Value = 507.42856
$ws.Range("K14").Value = 1522.28568
$ws.Range("M14").Value = -1349.28568

# --- Sheet GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H132").Value = 5151.4546
$ws.Range("I132").Value = 4074.6667
$ws.Range("K132").Value = 12224.0001
$ws.Range("M132").Value = -9694.000100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H55").Value = 3372.2727
$ws.Range("I55").Value = 2442.7144
$ws.Range("J55").Value = 4999
$ws.Range("K55").Value = 2442.7144
$ws.Range("L55").Value = 4999
$ws.Range("M55").Value = -2269.7144
$ws.Range("N55").Value = -5345

# --- Sheet WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Range("H13").Value = 1250
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = -360
$ws.Range("N13").Value = -2280
$ws.Range("H46").Value = 99995
$ws.Range("J46").Value = 99995
$ws.Range("L46").Value = 99995
$ws.Range("N46").Value = -100457
$ws.Range("H130").Value = 33750
$ws.Range("J130").Value = 32500
$ws.Range("L130").Value = 32500
$ws.Range("N130").Value = -42540
$ws.Range("H134").Value = 99995
$ws.Range("J134").Value = 99995
$ws.Range("L134").Value = 299985
$ws.Range("N134").Value = -305055
